# "harmonized similar tags to be the same"
# The #TAGS list block on the SwateTemplateMetadata sheet is extended: the
# single "plant growth protocol" tag (with its DPBO/NCIT term refs) is
# replaced by two harmonized tags - "Plant" (NCIT:C14258) and "growth"
# (GO:0040007) - kept alongside the existing "study" tag (NCIT:C63536,
# shortened from its full PURL form). The now-unused Term Source REF
# values in row 14 are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")

# Row 13: Tags Term Accession Number - shorten the existing "study" term ref
# to CURIE notation, and add the two new harmonized term refs.
$ws.Range("B13").Value = "NCIT:C63536"
$ws.Range("D12").Value = "growth"
$ws.Range("D13").Value = "GO:0040007"
$ws.Range("C12").Value = "Plant"
$ws.Range("C13").Value = "NCIT:C14258"

# Row 14: Tags Term Source REF (no longer needed - cleared)
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = ""

# Copy formatting from the existing column-B/C cells onto the newly used
# column D cells so the new cells match the rest of the block's style.
$ws.Range("B13").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4122) | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
